$d = $word.ActiveDocument
$sel = $word.Selection

$sel.TypeText('AI Impact Read-Me (markdown) file. Draft 1 (Q & A style) – one day looking back on this, we want to know why we did what we and more importantly how we did it. ')
$sel.TypeParagraph()
$sel.TypeText('Something to look back to as reference, to remember and reflect upon!')
$sel.TypeParagraph()
$sel.TypeText('The Team – Group B (“Hackathon One”)')
$sel.TypeParagraph()
$sel.TypeText('Who was in our group and what were their skillset/specialities?')
$sel.TypeParagraph()
$sel.TypeText('Our group members were assigned by our teacher and coach at Code Institute, John A. However, he left it up to us to split the group up into specific roles. We had different people from different backgrounds. There were 2 dashboard specialists, 2 with project management experience and one business analyst (me) but we all had the same data analyst knowledge that we learnt from being on the course.')
$sel.TypeParagraph()
$sel.TypeText('<Add names and roles here> ')
$sel.TypeParagraph()
$sel.TypeText('How did we approach task assignment?')
$sel.TypeParagraph()
$sel.TypeText('We worked out the best way to tackle this is to ask everyone to ‘muck in’ in two ways – a primary role and a secondary one. Primary is essentially their ‘bread and butter’ i.e. the skills they already possess, e.g. project management skills would project manage the team. The secondary would be a learning opportunity, so those say who wanted to learn a part of data analytics they were still unsure about (for me, I opted to work on the data pipeline, which to-date I’ve not yet used before).')
$sel.TypeParagraph()
$sel.TypeText('We had discussions about this and decided it would be the best thing to do and the fairest way to assign tasks. Afterall, this was a hackathon where time and deadlines were paramount, so it was really an opportunity to gain experience first and foremost, where learning would be a secondary goal. This would be fair enough suggestion, as our first hackathon – and in the end this was what we did. ')
$sel.TypeParagraph()
$sel.TypeText(' Tools and Applications')
$sel.TypeParagraph()
$sel.TypeText('What tools did we use?')
$sel.TypeParagraph()
$sel.TypeText(' For the most part, we used the tools provided to us as part of the data analytics course (See below table).')
$sel.TypeParagraph()
$sel.TypeText('How did we decide to use these tools?')
$sel.TypeParagraph()
$sel.TypeText('As this was predominantly a Python-based course, we chose to code in Python. The hackathon end-goal was to create a dashboard and at this point in our course we knew 3 dashboard applications (Streamlit, Power BI and Tableau).')
$sel.TypeParagraph()
$sel.TypeText('The Dataset')
$sel.TypeParagraph()
$sel.TypeText('Where did we find it? We found it on Kaggle. Details below.')
$sel.TypeParagraph()
$sel.TypeText('Why did we choose this dataset? Was it a good choice? How did you make it fit-for-purpose?')
$sel.TypeParagraph()
$sel.TypeText('With what’s going on in the modern world, as well as what the course was teaching us, we thought it would be a good opportunity to see if we were able to work out if we still had jobs in 2030 or would it have already been taken by AI.')
$sel.TypeParagraph()
$sel.TypeText('What did we learn from the dataset?')
$sel.TypeParagraph()
$sel.TypeText('Initially, when we chose the dataset, it looked promising – as it was to do with the subject we were interested in and currently learning about. However, after we had a closer look at the data, we found it to be ‘too synthetic’ and did not accurately represent what we were trying to predict.')
$sel.TypeParagraph()
$sel.TypeText('Did it matter that it was synthesized data?')
$sel.TypeParagraph()
$sel.TypeText('The problem wasn’t that it was synthetic data. As from my previous project looking at synth loans data – these can be quite accurate. Here we found the problem was that the average of ALL numerical fields was the same. No matter which field you chose (details will follow later), there was no way to distinguish say a doctor’s wages against those of a nurse as the data suggested the average pay for the two job titles were one and the same. So, from this, we could very well conclude that a doctor would have the same chance as a nurse to lose their job to AI based on their pay scales if their average pay was the same.   ')
$sel.TypeParagraph()
$sel.TypeText('What was the quality of the data like? How did we clean it?')
$sel.TypeParagraph()
$sel.TypeText('Data quality was fine with minimum cleaning routines applied (basic checks for missing data showed that there wasn’t any).')
$sel.TypeParagraph()
$sel.TypeParagraph()
$sel.TypeText('The Opportunity')
$sel.TypeParagraph()
$sel.TypeText('From our brainstorming sessions, what questions did we ask?')
$sel.TypeParagraph()
$sel.TypeText('[Insert questions here]')
$sel.TypeParagraph()
$sel.TypeText('The “Quick & Dirty” (Primary/Initial) Analysis ')
$sel.TypeParagraph()
$sel.TypeText('How did we approach our analysis?')
$sel.TypeParagraph()
$sel.TypeText('Used Microsoft Excel to see the dataset in its entirety. This approach also showed us, by using simple line graphs) which jobs were at high risk and which jobs were safe in 2030. It also showed us that the dataset was a full set without any missing or strange values we had to deal with.')
$sel.TypeParagraph()
$sel.TypeText('At first glance, what does it tell us about AI?')
$sel.TypeParagraph()
$sel.TypeText('What hurdles did we have to overcome?')
$sel.TypeParagraph()
$sel.TypeText('Was the data able to answer the questions from your brainstorming sessions?')
$sel.TypeParagraph()
$sel.TypeParagraph()
$sel.TypeText('The Analysis (Secondary) – The Deep Dive using Python & its Data Libraries')
$sel.TypeParagraph()
$sel.TypeText('How did we refine our analysis?')
$sel.TypeParagraph()
$sel.TypeText('What hurdles did we have to overcome?')
$sel.TypeParagraph()
$sel.TypeText('Where there any new features added? Why were they added? What was the impact?')
$sel.TypeParagraph()
$sel.TypeText('Dashboard Coding & Design')
$sel.TypeParagraph()
$sel.TypeParagraph()
$sel.TypeText('How did your group come up with the ideas for the dashboard layout?')
$sel.TypeParagraph()
$sel.TypeText('Why did you choose two applications? What was your reasoning?')
$sel.TypeParagraph()
$sel.TypeParagraph()
$sel.TypeText('The Presentation')
$sel.TypeParagraph()
$sel.TypeText('How did we decide to present the data – on the dashboard and why did we use two applications?')
$sel.TypeParagraph()
$sel.TypeParagraph()
$sel.TypeText('Next Steps')
$sel.TypeParagraph()
$sel.TypeText('What’s Next?')
$sel.TypeParagraph()

$bothIdx = 1,2,5,6,8,9,10,14,18,22
foreach ($idx in $bothIdx) {
    $d.Paragraphs($idx).Range.ParagraphFormat.Alignment = 3
}
